$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column W (day 22) values for the Bibi stores, and the corresponding
# row totals in column AG which sum the daily values.
$ws.Range("W2").Value = 6826.2
$ws.Range("AG2").Value = 191499.76

$ws.Range("W3").Value = 5892
$ws.Range("AG3").Value = 96455.00999999999

$ws.Range("W4").Value = 1815
$ws.Range("AG4").Value = 66225.89999999999

$ws.Range("W5").Value = 1812
$ws.Range("AG5").Value = 59325.79

$ws.Range("W6").Value = 16345.2
$ws.Range("AG6").Value = 413506.46
